
$d = $word.ActiveDocument
$s = $d.Styles.Add("NumTest", 4)
Write-Output ("Visibility before=" + $s.Visibility)
$s.Visibility = $true
Write-Output ("Visibility after true=" + $s.Visibility)
